$d = $word.ActiveDocument

# Locate the inline picture that was inserted as a "call tree" image
# (descr = "A group of text boxes\n\nDescription automatically generated")
# and remove it, along with the run that hosts it (the w:lastRenderedPageBreak
# + w:drawing run). The following w:br w:type="page" run in the same
# paragraph must be left untouched.
for ($i = $d.InlineShapes.Count; $i -ge 1; $i--) {
    $shape = $d.InlineShapes.Item($i)
    if ($shape.Type -eq 3 -and $shape.AlternativeText -like "*group of text boxes*") {
        $shape.Delete()
    }
}
